$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list: update each coin's 'Price' (column D) and
# 'Volume(1h)' (column E) text values to the latest scraped figures.
#
# A few Price values are plain decimals (e.g. 561.58). Excel's COM layer
# auto-converts such text into a Number when assigned directly, which
# would strip significant trailing zeros and introduce floating-point
# noise. To keep them as plain text (as in the source sheet), we briefly
# mark the cell as Text before assigning, then restore the default
# 'Normal' style so no visible formatting change is left behind.

$ws.Range("D2").Value = "59.296.37"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Value = "2.993.05"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.88%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "2.980.37"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.14"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.80%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.88%  "
$ws.Range("E13").Value = "  +1.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("E15").Value = "  +1.76%  "
$ws.Range("D16").Value = "3.491.72"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.03%  "
$ws.Range("D18").Value = "2.995.35"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "59.284.95"
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.77%  "
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -7.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.992"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("D36").Value = "0.0₃0761"
$ws.Range("E36").Value = "  +9.28%  "
$ws.Range("E37").Value = "  -2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.66%  "
$ws.Range("E39").Value = "  +2.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "401.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0350"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.69%  "
$ws.Range("D43").Value = "2.754.16"
$ws.Range("E43").Value = "  +3.15%  "
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +17.24%  "
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
